$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: D0.1596976872874879 -> D0.18565903795358007
$ws.Range("A2").Value = "D0.18565903795358007"

# D2:D4: Sat, 10 Dec 2022 18:18:45 +0530 -> Mon, 19 Dec 2022 05:44:49 -0800
$ws.Range("D2").Value = "Mon, 19 Dec 2022 05:44:49 -0800"
$ws.Range("D3").Value = "Mon, 19 Dec 2022 05:44:49 -0800"
$ws.Range("D4").Value = "Mon, 19 Dec 2022 05:44:49 -0800"

# A3: D0.7566317244104037 -> D0.3123609309210864
$ws.Range("A3").Value = "D0.3123609309210864"

# A4: D0.5185698880184182 -> D0.6471154244116394
$ws.Range("A4").Value = "D0.6471154244116394"
